$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("company_list")

# Row 2
$ws.Range("D2").Value = 131536
$ws.Range("E2").Value = 5830
$ws.Range("F2").Value = 5830
$ws.Range("G2").Value = 4238
$ws.Range("H2").Value = 2919
$ws.Range("I2").Value = 2900
$ws.Range("J2").Value = 19
$ws.Range("K2").Value = 138275
$ws.Range("L2").Value = 67677
$ws.Range("M2").Value = 70597
$ws.Range("N2").Value = 69078
$ws.Range("O2").Value = 1519
$ws.Range("P2").Value = 1394
$ws.Range("Q2").Value = 7010
$ws.Range("R2").Value = -9578
$ws.Range("S2").Value = 2505
$ws.Range("T2").Value = 9591
$ws.Range("U2").Value = -2581
$ws.Range("V2").Value = 38397
$ws.Range("W2").Value = 4.43
$ws.Range("X2").Value = 2.22
$ws.Range("Y2").Value = 4.3
$ws.Range("Z2").Value = 2.17
$ws.Range("AA2").Value = 95.86
$ws.Range("AB2").Value = 3937.28
$ws.Range("AC2").Value = 10404
$ws.Range("AD2").Value = 19.51
$ws.Range("AE2").Value = 247910
$ws.Range("AF2").Value = 0.82
$ws.Range("AG2").Value = 1500
$ws.Range("AH2").Value = 0.74
$ws.Range("AI2").Value = 14.41
$ws.Range("AJ2").Value = 27875819

# Row 3
$ws.Range("D3").Value = 136400
$ws.Range("E3").Value = 5038
$ws.Range("F3").Value = 5038
$ws.Range("G3").Value = 6938
$ws.Range("H3").Value = 4559
$ws.Range("I3").Value = 4547
$ws.Range("J3").Value = 12
$ws.Range("K3").Value = 144938
$ws.Range("L3").Value = 72544
$ws.Range("M3").Value = 72394
$ws.Range("N3").Value = 69861
$ws.Range("O3").Value = 2533
$ws.Range("P3").Value = 1394
$ws.Range("Q3").Value = 7339
$ws.Range("R3").Value = -10039
$ws.Range("S3").Value = 2791
$ws.Range("T3").Value = 10064
$ws.Range("U3").Value = -2725
$ws.Range("V3").Value = 41028
$ws.Range("W3").Value = 3.69
$ws.Range("X3").Value = 3.34
$ws.Range("Y3").Value = 6.54
$ws.Range("Z3").Value = 3.22
$ws.Range("AA3").Value = 100.21
$ws.Range("AB3").Value = 4226.05
$ws.Range("AC3").Value = 16312
$ws.Range("AD3").Value = 11.59
$ws.Range("AE3").Value = 250718
$ws.Range("AF3").Value = 0.75
$ws.Range("AG3").Value = 1500
$ws.Range("AH3").Value = 0.79
$ws.Range("AI3").Value = 9.19
$ws.Range("AJ3").Value = 27875819

# Row 4
$ws.Range("D4").Value = 146151
$ws.Range("E4").Value = 5686
$ws.Range("F4").Value = 5469
$ws.Range("G4").Value = 5068
$ws.Range("H4").Value = 3816
$ws.Range("I4").Value = 3762
$ws.Range("J4").Value = 54
$ws.Range("K4").Value = 154301
$ws.Range("L4").Value = 73059
$ws.Range("M4").Value = 81242
$ws.Range("N4").Value = 76962
$ws.Range("O4").Value = 4280
$ws.Range("P4").Value = 1394
$ws.Range("Q4").Value = 7356
$ws.Range("R4").Value = -9161
$ws.Range("S4").Value = 1846
$ws.Range("T4").Value = 6002
$ws.Range("U4").Value = 1354
$ws.Range("V4").Value = 38762
$ws.Range("W4").Value = 3.89
$ws.Range("X4").Value = 2.61
$ws.Range("Y4").Value = 5.13
$ws.Range("Z4").Value = 2.55
$ws.Range("AA4").Value = 89.93000000000001
$ws.Range("AB4").Value = 4444.69
$ws.Range("AC4").Value = 13497
$ws.Range("AD4").Value = 13.56
$ws.Range("AE4").Value = 276203
$ws.Range("AF4").Value = 0.66
$ws.Range("AG4").Value = 1500
$ws.Range("AH4").Value = 0.82
$ws.Range("AI4").Value = 11.11
$ws.Range("AJ4").Value = 27875819

# Row 5
$ws.Range("D5").Value = 155149
$ws.Range("E5").Value = 5849
$ws.Range("F5").Value = 5849
$ws.Range("G5").Value = 7997
$ws.Range("H5").Value = 6279
$ws.Range("I5").Value = 6161
$ws.Range("J5").Value = 119
$ws.Range("K5").Value = 160665
$ws.Range("L5").Value = 72951
$ws.Range("M5").Value = 87714
$ws.Range("N5").Value = 82417
$ws.Range("O5").Value = 5297
$ws.Range("P5").Value = 1394
$ws.Range("Q5").Value = 8982
$ws.Range("R5").Value = -6218
$ws.Range("S5").Value = -1214
$ws.Range("T5").Value = 7228
$ws.Range("U5").Value = 1754
$ws.Range("V5").Value = 36593
$ws.Range("W5").Value = 3.77
$ws.Range("X5").Value = 4.05
$ws.Range("Y5").Value = 7.73
$ws.Range("Z5").Value = 3.99
$ws.Range("AA5").Value = 83.17
$ws.Range("AB5").Value = 4808.41
$ws.Range("AC5").Value = 22101
$ws.Range("AD5").Value = 12.26
$ws.Range("AE5").Value = 295780
$ws.Range("AF5").Value = 0.92
$ws.Range("AG5").Value = 1750
$ws.Range("AH5").Value = 0.65
$ws.Range("AI5").Value = 7.91
$ws.Range("AJ5").Value = 27875819

# Row 6
$ws.Range("D6").Value = 170491
$ws.Range("E6").Value = 4628
$ws.Range("F6").Value = 4628
$ws.Range("G6").Value = 5850
$ws.Range("H6").Value = 4762
$ws.Range("I6").Value = 4502
$ws.Range("K6").Value = 167539
$ws.Range("L6").Value = 78964
$ws.Range("M6").Value = 88575
$ws.Range("N6").Value = 81723
$ws.Range("P6").Value = 1394
$ws.Range("Q6").Value = 7698
$ws.Range("R6").Value = -8168
$ws.Range("S6").Value = 1028
$ws.Range("T6").Value = 8946
$ws.Range("U6").Value = -1248
$ws.Range("V6").Value = 38186
$ws.Range("W6").Value = 2.71
$ws.Range("X6").Value = 2.79
$ws.Range("Y6").Value = 5.49
$ws.Range("Z6").Value = 2.9
$ws.Range("AA6").Value = 89.15000000000001
$ws.Range("AB6").Value = 5029.24
$ws.Range("AC6").Value = 16150
$ws.Range("AD6").Value = 11.3
$ws.Range("AE6").Value = 293291
$ws.Range("AF6").Value = 0.62
$ws.Range("AG6").Value = 2000
$ws.Range("AH6").Value = 1.1
$ws.Range("AI6").Value = 12.38
$ws.Range("AJ6").Value = 27875819

# Row 7
$ws.Range("D7").Value = 189324
$ws.Range("E7").Value = 2039
$ws.Range("G7").Value = 2358
$ws.Range("H7").Value = 1949
$ws.Range("I7").Value = 1878
$ws.Range("K7").Value = 196901
$ws.Range("L7").Value = 98299
$ws.Range("M7").Value = 98602
$ws.Range("N7").Value = 86886
$ws.Range("P7").Value = 1391
$ws.Range("Q7").Value = 11071
$ws.Range("R7").Value = -13279
$ws.Range("S7").Value = 5076
$ws.Range("T7").Value = 7748
$ws.Range("U7").Value = 1665
$ws.Range("W7").Value = 1.08
$ws.Range("X7").Value = 1.03
$ws.Range("Y7").Value = 2.23
$ws.Range("Z7").Value = 1.07
$ws.Range("AA7").Value = 99.69
$ws.Range("AC7").Value = 6737
$ws.Range("AD7").Value = 16.48
$ws.Range("AE7").Value = 322226
$ws.Range("AF7").Value = 0.34
$ws.Range("AG7").Value = 1858
$ws.Range("AH7").Value = 1.67
$ws.Range("AI7").Value = 27.58

# Row 8
$ws.Range("D8").Value = 202156
$ws.Range("E8").Value = 2851
$ws.Range("G8").Value = 3246
$ws.Range("H8").Value = 2542
$ws.Range("I8").Value = 2422
$ws.Range("K8").Value = 199534
$ws.Range("L8").Value = 99097
$ws.Range("M8").Value = 100438
$ws.Range("N8").Value = 88591
$ws.Range("P8").Value = 1391
$ws.Range("Q8").Value = 11118
$ws.Range("R8").Value = -7714
$ws.Range("S8").Value = -2001
$ws.Range("T8").Value = 6520
$ws.Range("U8").Value = 3637
$ws.Range("W8").Value = 1.41
$ws.Range("X8").Value = 1.26
$ws.Range("Y8").Value = 2.76
$ws.Range("Z8").Value = 1.28
$ws.Range("AA8").Value = 98.66
$ws.Range("AC8").Value = 8688
$ws.Range("AD8").Value = 12.78
$ws.Range("AE8").Value = 328551
$ws.Range("AF8").Value = 0.34
$ws.Range("AG8").Value = 1903
$ws.Range("AH8").Value = 1.71
$ws.Range("AI8").Value = 21.9

# Row 9
$ws.Range("D9").Value = 214154
$ws.Range("E9").Value = 3602
$ws.Range("G9").Value = 4051
$ws.Range("H9").Value = 3149
$ws.Range("I9").Value = 3029
$ws.Range("K9").Value = 203997
$ws.Range("L9").Value = 101294
$ws.Range("M9").Value = 102703
$ws.Range("N9").Value = 90721
$ws.Range("P9").Value = 1391
$ws.Range("Q9").Value = 11483
$ws.Range("R9").Value = -7923
$ws.Range("S9").Value = -971
$ws.Range("T9").Value = 6656
$ws.Range("U9").Value = 3811
$ws.Range("W9").Value = 1.68
$ws.Range("X9").Value = 1.47
$ws.Range("Y9").Value = 3.38
$ws.Range("Z9").Value = 1.56
$ws.Range("AA9").Value = 98.63
$ws.Range("AC9").Value = 10868
$ws.Range("AD9").Value = 10.21
$ws.Range("AE9").Value = 336451
$ws.Range("AF9").Value = 0.33
$ws.Range("AG9").Value = 1931
$ws.Range("AH9").Value = 1.74
$ws.Range("AI9").Value = 17.76
